$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated odds/values per diff (row 2..13, various columns)
$ws.Range("H2").Value = 2.38
$ws.Range("L2").Value = 1.4
$ws.Range("Q2").Value = 1.93
$ws.Range("S2").Value = 3.35
$ws.Range("T2").Value = 1.72
$ws.Range("U2").Value = 2.1
$ws.Range("AH2").Value = 21
$ws.Range("AM2").Value = 110
$ws.Range("G3").Value = 2.56
$ws.Range("H3").Value = 2.8
$ws.Range("J3").Value = 3.35
$ws.Range("N3").Value = 3.45
$ws.Range("P3").Value = 1.86
$ws.Range("Q3").Value = 1.74
$ws.Range("S3").Value = 2.84
$ws.Range("W3").Value = 1.64
$ws.Range("N4").Value = 4
$ws.Range("F5").Value = 8.4
$ws.Range("G5").Value = 8.6
$ws.Range("N5").Value = 4.9
$ws.Range("O5").Value = 1.24
$ws.Range("S5").Value = 2.9
$ws.Range("W5").Value = 1.13
$ws.Range("Y5").Value = 9
$ws.Range("AG5").Value = 30
$ws.Range("AJ5").Value = 280
$ws.Range("P6").Value = 2.5
$ws.Range("S6").Value = 2.4
$ws.Range("I7").Value = 1.9
$ws.Range("L7").Value = 1.34
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 4.7
$ws.Range("O7").Value = 1.25
$ws.Range("P7").Value = 2.26
$ws.Range("Q7").Value = 1.77
$ws.Range("R7").Value = 1.5
$ws.Range("S7").Value = 2.94
$ws.Range("T7").Value = 1.73
$ws.Range("U7").Value = 2.3
$ws.Range("X7").Value = 18
$ws.Range("Z7").Value = 12
$ws.Range("AC7").Value = 8.800000000000001
$ws.Range("AF7").Value = 34
$ws.Range("AG7").Value = 16.5
$ws.Range("AI7").Value = 30
$ws.Range("AK7").Value = 50
$ws.Range("AL7").Value = 55
$ws.Range("AN7").Value = 46
$ws.Range("J8").Value = 1.03
$ws.Range("N8").Value = 1.1
$ws.Range("Q8").Value = 1.27
$ws.Range("S8").Value = 1.27
$ws.Range("T8").Value = 1.05
$ws.Range("U8").Value = 1.05
$ws.Range("F9").Value = 2.36
$ws.Range("G9").Value = 2.38
$ws.Range("N9").Value = 3.7
$ws.Range("O9").Value = 1.35
$ws.Range("P9").Value = 1.9
$ws.Range("Q9").Value = 2.08
$ws.Range("R9").Value = 1.35
$ws.Range("W9").Value = 1.72
$ws.Range("X9").Value = 12.5
$ws.Range("Y9").Value = 13
$ws.Range("S10").Value = 2.9
$ws.Range("X10").Value = 21
$ws.Range("AA10").Value = 9.199999999999999
$ws.Range("AE10").Value = 14.5
$ws.Range("F11").Value = 2.72
$ws.Range("H11").Value = 2.86
$ws.Range("S11").Value = 3.45
$ws.Range("X11").Value = 13.5
$ws.Range("G12").Value = 4.2
$ws.Range("H12").Value = 1.97
$ws.Range("I12").Value = 1.99
$ws.Range("J12").Value = 4
$ws.Range("K12").Value = 4.1
$ws.Range("P12").Value = 2.36
$ws.Range("R12").Value = 1.53
$ws.Range("V12").Value = 2.02
$ws.Range("W12").Value = 1.31
$ws.Range("AC12").Value = 9
$ws.Range("G13").Value = 5
$ws.Range("I13").Value = 1.8
$ws.Range("J13").Value = 4.1
$ws.Range("N13").Value = 5.8
$ws.Range("O13").Value = 1.2
$ws.Range("R13").Value = 1.66
$ws.Range("S13").Value = 2.44
$ws.Range("V13").Value = 2.24
$ws.Range("X13").Value = 23
$ws.Range("AB13").Value = 24
$ws.Range("AH13").Value = 15.5
$ws.Range("AJ13").Value = 100
$ws.Range("AL13").Value = 48
